$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose values look numeric,
# so Excel keeps them as exact text instead of converting to a number.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value = '27.950.55'
$ws.Range('D3').Value = '1.767.22'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '328.47'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '0.4686'
$ws.Range('E7').Value = '  +2.00%  '
$ws.Range('D8').Value = '0.3529'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('D9').Value = '43.72'
$ws.Range('E9').Value = '  +4.64%  '
$ws.Range('D10').Value = '0.07382'
$ws.Range('E10').Value = '  -1.38%  '
$ws.Range('D11').Value = '1.083'
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('D12').Value = '0.9996'
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').Value = '20.62'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = '6.005'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').Value = '7.179'
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').Value = '1.764.49'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '92.25'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '0.00001054'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = '0.06422'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '0.9999'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '16.92'
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('D22').Value = '5.780'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '27.978.54'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').Value = '11.13'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').Value = '2.156'
$ws.Range('E25').Value = '  +3.46%  '
$ws.Range('D26').Value = '162.90'
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('D27').Value = '20.01'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('D28').Value = '1.964.53'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').Value = '2.180'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = '122.82'
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('D31').Value = '1.073'
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('D32').Value = '0.09303'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').Value = '3.646'
$ws.Range('E33').Value = '  -0.64%  '
$ws.Range('D34').Value = '5.548'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('D35').Value = '11.69'
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('D36').Value = '0.02267'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').Value = '0.06092'
$ws.Range('E37').Value = '  -1.67%  '
$ws.Range('D38').Value = '0.2067'
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('D39').Value = '4.912'
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.6150'
$ws.Range('E40').Value = '  -2.48%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.187'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').Value = '1.427'
$ws.Range('E42').Value = '  +2.70%  '
$ws.Range('D43').Value = '7.760'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').Value = '13.08'
$ws.Range('E44').Value = '  -1.54%  '
$ws.Range('D45').Value = '3.740'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Value = '0.5792'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D47').Value = '123.59'
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('D48').Value = '1.931'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '1.125'
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.06812'
$ws.Range('E50').Value = '  -1.57%  '
$ws.Range('D51').Value = '72.14'
$ws.Range('E51').Value = '  +0.04%  '
